$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.297606706619263
$ws.Range("B1").Value = 2.63785457611084
$ws.Range("C1").Value = 1.560036301612854
$ws.Range("D1").Value = 1.262665271759033
$ws.Range("E1").Value = 1.174902677536011
